# Add a new "{SenderDepartment}" placeholder paragraph right after the
# existing "{SenderName}" paragraph at the top of the letterhead block
# (mirrors the Arial rFonts formatting already used by the surrounding
# sender-address paragraphs).

$d = $word.ActiveDocument

# Locate the paragraph that holds the {SenderName} placeholder.
$range = $d.Content
$found = $range.Find.Execute("{SenderName}", $false, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)

$senderNamePara = $range.Paragraphs.First
$endOfSenderNamePara = $senderNamePara.Range.End

# Insert a brand-new paragraph immediately after it; Word clones the
# paragraph/run formatting (Arial rFonts) from the paragraph it split off.
$insertPoint = $d.Range($endOfSenderNamePara, $endOfSenderNamePara)
$insertPoint.InsertParagraphAfter()

# Fill the freshly created (empty) paragraph with the department placeholder.
$newParaStart = $endOfSenderNamePara + 1
$newParaRange = $d.Range($newParaStart, $newParaStart)
$newParaRange.InsertAfter("{SenderDepartment}")
